# Tradução da história principal
# Fill in the Spanish/Portuguese translation column (C) of the message
# table: header C1 "TRADUÇÃO" -> "TRADUCCIÓN", and each row's translated
# text in column C (column B, the ORIGINAL text, stays untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Cells.Item(1, 3).Value = "TRADUCCIÓN"

# row -> translated text for column C (rows keyed by their OFFSET in col A)
$translations = @{
    7258 = "Traje recebido <Color:8>Tatuagem do Dragão<Color:Default>.`n\n"
    7313 = "Recebido <Color:8>Toughness Emperor<Color:Default>.`n\n"
    7366 = "Recebido <Color:8>Tauriner ++<Color:Default>."
    7413 = "Recebido <Color:8>Staminan Royale<Color:Default>."
    7464 = "Recebido <Color:8>Modified Model Gun<Color:Default>."
    7518 = "Recebido <Color:8>Prato de Ouro<Color:Default>."
    7565 = "Traje recebido <Color:8>Prisioneiro Fugitivo<Color:Default>."
    7621 = "Recebido <Color:8>Extra Balanced Motor<Color:Default>."
    7677 = "Recebido <Color:8>Extra Slim Tires<Color:Default>."
    7729 = "Recebido <Color:8>Boost Gears<Color:Default>."
    7776 = "Recebido <Color:8>Speed Frame Plus<Color:Default>."
    7828 = "Recebido <Color:8>Stone of Enduring<Color:Default>."
    7882 = "Traje recebido <Color:8>Terno Preto<Color:Default>."
    7934 = "Recebido <Color:8>Toughness Infinity<Color:Default>."
    7988 = "Recebido <Color:8>Tauriner Maximum<Color:Default>."
    8040 = "Recebido <Color:8>Staminan Spark<Color:Default>."
    8090 = "Recebido <Color:8>Canhão de Peixe Espada<Color:Default>."
    8139 = "Recebido <Color:8>Prato de Platina<Color:Default>."
    8190 = "Traje recebido <Color:8>Casaco de Pele de Cobra<Color:Default>."
    8248 = "Recebido <Color:8>Killer Bee<Color:Default>."
    8294 = "Recebido <Color:8>Godspeed Motor<Color:Default>."
    8344 = "Recebido <Color:8>Super Slim Tires<Color:Default>."
    8396 = "Recebido <Color:8>Godspeed Gears Plus<Color:Default>."
    8451 = "Recebido <Color:8>New Bumper Plate<Color:Default>."
    8503 = "Recebido <Color:8>Pedra do Sacrifício<Color:Default>."
    8555 = "Recebido <Color:8>Calming Towel<Color:Default>."
    8605 = "Recebido <Color:8>Colar Magnético<Color:Default>."
    8659 = "Recebido <Color:8>Amuleto da Família Dojima<Color:Default>.`n\n"
    8716 = "Recebido <Color:8>Talismã do Deus da Guerra<Color:Default>."
}

$dim = $ws.UsedRange
$rowCount = $dim.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $offset = $ws.Cells.Item($r, 1).Value2
    if ($null -eq $offset) { continue }
    $key = [int]$offset
    if ($translations.ContainsKey($key)) {
        $ws.Cells.Item($r, 3).Value = $translations[$key]
    }
}
